# "forgot to update sprint backlog"
# Day 4 (F column) task-hours for "Task #2: documentation of functions" (row 8)
# was never entered - fill in the missing 2 hours, matching the red
# "value entered" formatting used elsewhere in the sheet (e.g. F11, D7, E7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$f8 = $ws.Range("F8")
$f8.Value = 2
$f8.Font.Color = 255   # RGB(255,0,0) - red, same as the other entered-value cells

# Re-enter the "Task Sum" row formula across C13:L13 (the daily totals), which
# recalculates F13 (and, in turn, N13/B15 "Total Work") to account for the new
# Day 4 hours above.
$ws.Range("C13:L13").Formula = "=SUM(C3:C10)"

# Leave the selection where the edit was made.
$ws.Range("F8").Select()
